# Updates the cryptos price/volume table to reflect the latest GitHub
# Actions scrape. Columns: A=rank(idx) B=Coin C=Link D=Price E=Volume(1h)
#
# Column D ("Price") values are stored as TEXT in the workbook (e.g.
# "26.140.99", "1.005") even though many of them look like plain decimal
# numbers. Excel's COM layer auto-coerces a numeric-looking string typed
# into Range.Value into an actual number (losing the literal text and
# adding a float-precision artifact), so those assignments are wrapped in
# a force-text helper: mark the cell as Text ("@") before writing, then
# restore the "Normal" style afterwards so no stray number format sticks
# around on the cell (matches the original workbook, where none of these
# cells carry an explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-Cell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# Row 2 - Bitcoin
Set-TextCell 2 4 "26.140.99"
Set-Cell     2 5 "  -1.17%  "

# Row 3 - Ethereum
Set-TextCell 3 4 "1.657.83"
Set-Cell     3 5 "  -1.07%  "

# Row 4 - TetherUSD
Set-TextCell 4 4 "1.005"
Set-Cell     4 5 "  +0.32%  "

# Row 5 - BNB
Set-TextCell 5 4 "216.12"
Set-Cell     5 5 "  -1.49%  "

# Row 6 - XRP
Set-TextCell 6 4 "0.5199"
Set-Cell     6 5 "  -2.10%  "

# Row 7 - USDC
Set-TextCell 7 4 "1.004"

# Row 8 - Cardano
Set-TextCell 8 4 "0.2626"
Set-Cell     8 5 "  -2.66%  "

# Row 9 - Dogecoin
Set-TextCell 9 4 "0.06264"
Set-Cell     9 5 "  -2.04%  "

# Row 10 - Solana
Set-TextCell 10 4 "20.71"
Set-Cell     10 5 "  -5.01%  "

# Row 11 - TRON
Set-TextCell 11 4 "0.07719"
Set-Cell     11 5 "  -1.07%  "

# Row 12 - now Polkadot (was WrappedEther)
Set-Cell     12 2 "Polkadot"
Set-Cell     12 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 12 4 "4.423"
Set-Cell     12 5 "  -1.88%  "

# Row 13 - now WrappedEther (was Polkadot)
Set-Cell     13 2 "WrappedEther"
Set-Cell     13 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell 13 4 "1.642.89"
Set-Cell     13 5 "  -1.99%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextCell 14 4 "1.885.48"
Set-Cell     14 5 "  -1.04%  "

# Row 15 - Polygon
Set-TextCell 15 4 "0.5407"
Set-Cell     15 5 "  -3.14%  "

# Row 16 - ShibaInu
Set-TextCell 16 4 "0.0₅8126"
Set-Cell     16 5 "  -2.47%  "

# Row 17 - Litecoin
Set-TextCell 17 4 "64.56"
Set-Cell     17 5 "  -1.63%  "

# Row 18 - WrappedBTC
Set-TextCell 18 4 "26.177.38"
Set-Cell     18 5 "  -1.15%  "

# Row 19 - Dai
Set-Cell 19 5 "  +0.35%  "

# Row 20 - Uniswap
Set-TextCell 20 4 "4.611"
Set-Cell     20 5 "  -3.60%  "

# Row 21 - BitcoinCash
Set-TextCell 21 4 "191.30"
Set-Cell     21 5 "  -0.83%  "

# Row 22 - Avalanche
Set-TextCell 22 4 "10.04"
Set-Cell     22 5 "  -2.38%  "

# Row 23 - Chainlink
Set-TextCell 23 4 "6.036"
Set-Cell     23 5 "  -4.51%  "

# Row 24 - BinanceUSD
Set-TextCell 24 4 "1.006"
Set-Cell     24 5 "  +0.39%  "

# Row 25 - Monero
Set-TextCell 25 4 "139.31"
Set-Cell     25 5 "  -0.50%  "

# Row 26 - Stellar
Set-TextCell 26 4 "0.1226"

# Row 27 - Cosmos
Set-Cell 27 5 "  -2.99%  "

# Row 28 - EthereumClassic
Set-Cell 28 5 "  -1.52%  "

# Row 29 - Toncoin
Set-TextCell 29 4 "1.400"
Set-Cell     29 5 "  -3.02%  "

# Row 30 - Hedera
Set-TextCell 30 4 "0.05950"
Set-Cell     30 5 "  -4.96%  "

# Row 31 - PancakeSwap
Set-TextCell 31 4 "1.267"
Set-Cell     31 5 "  -1.28%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextCell 32 4 "3.540"
Set-Cell     32 5 "  -1.89%  "

# Row 33 - Filecoin
Set-TextCell 33 4 "3.248"
Set-Cell     33 5 "  -5.91%  "

# Row 34 - LidoDAOToken
Set-TextCell 34 4 "1.604"
Set-Cell     34 5 "  -5.16%  "

# Row 35 - ARBITRUM
Set-TextCell 35 4 "0.9633"
Set-Cell     35 5 "  -4.72%  "

# Row 36 - HuobiToken
Set-Cell 36 5 "  +0.11%  "

# Row 37 - MXToken
Set-TextCell 37 4 "2.774"
Set-Cell     37 5 "  -0.42%  "

# Row 38 - ImmutableX
Set-TextCell 38 4 "0.5651"
Set-Cell     38 5 "  -8.51%  "

# Row 39 - VeChain
Set-TextCell 39 4 "0.01593"
Set-Cell     39 5 "  -2.29%  "

# Row 40 - FraxShare
Set-TextCell 40 4 "5.967"
Set-Cell     40 5 "  -3.19%  "

# Row 41 - TrustWalletToken
Set-TextCell 41 4 "0.8555"
Set-Cell     41 5 "  -0.90%  "

# Row 42 - PaxDollar
Set-Cell 42 5 "  +0.34%  "

# Row 43 - Maker
Set-TextCell 43 4 "1.009.76"
Set-Cell     43 5 "  -7.85%  "

# Row 44 - Quant
Set-TextCell 44 4 "100.33"
Set-Cell     44 5 "  -0.21%  "

# Row 45 - RocketPoolETH
Set-TextCell 45 4 "1.799.99"

# Row 46 - BabyDogeCoin
Set-Cell 46 5 "  -0.20%  "

# Row 47 - Aave
Set-TextCell 47 4 "56.53"
Set-Cell     47 5 "  -3.61%  "

# Row 48 - Frax
Set-TextCell 48 4 "1.004"
Set-Cell     48 5 "  +0.06%  "

# Row 49 - EnergySwap
Set-TextCell 49 4 "7.967"
Set-Cell     49 5 "  -2.15%  "

# Row 50 - Cronos
Set-Cell 50 5 "  -0.48%  "

# Row 51 - Mantle
Set-TextCell 51 4 "0.4198"
Set-Cell     51 5 "  -0.87%  "
